$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "Done" status in column E for rows 5, 7, 8 and 9, matching the
# formatting (wrap text) already used by the other "Done" cells in the sheet.
$doneCells = @("E5", "E7", "E8", "E9")
foreach ($cellRef in $doneCells) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "Done"
    $rng.WrapText = $true
}

# Update the active selection to match the new position (E7) recorded in the sheet view.
$ws.Range("E7").Select()
